$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder/replace the header row:
# Before: Mesa | Producto | Cantidad | Precio | Fecha_Hora | Total | Estado | Categoría | Metodo_Pago | Referencia
# After:  Mesa | Producto | Cantidad | Precio | Categoría  | Fecha_Hora | Estado | Total

$ws.Range("E1").Value = "Categoría"
$ws.Range("F1").Value = "Fecha_Hora"
$ws.Range("G1").Value = "Estado"
$ws.Range("H1").Value = "Total"

# Remove the now-obsolete trailing columns (Metodo_Pago, Referencia) so the
# used range shrinks back down to A1:H1.
$ws.Columns.Item(10).Delete()
$ws.Columns.Item(9).Delete()
